$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(795, 1).Value = 794
$ws.Cells.Item(795, 2).Value = 96
$ws.Cells.Item(795, 3).Value = 43213
$ws.Cells.Item(795, 4).Value = 'April'
$ws.Cells.Item(795, 5).Value = 2018
$ws.Cells.Item(795, 6).Value = 'Monday'
$ws.Cells.Item(795, 7).Value = 'Shoulder Press'
$ws.Cells.Item(795, 8).Value = 25
$ws.Cells.Item(795, 9).Value = 4
$ws.Cells.Item(795, 10).Value = 8
$ws.Cells.Item(795, 11).Value = 'Shoulders'

$ws.Cells.Item(796, 1).Value = 795
$ws.Cells.Item(796, 2).Value = 96
$ws.Cells.Item(796, 3).Value = 43213
$ws.Cells.Item(796, 4).Value = 'April'
$ws.Cells.Item(796, 5).Value = 2018
$ws.Cells.Item(796, 6).Value = 'Monday'
$ws.Cells.Item(796, 7).Value = 'Shoulder Shrug'
$ws.Cells.Item(796, 8).Value = 25
$ws.Cells.Item(796, 9).Value = 4
$ws.Cells.Item(796, 10).Value = 8
$ws.Cells.Item(796, 11).Value = 'Shoulders'

$ws.Cells.Item(797, 1).Value = 796
$ws.Cells.Item(797, 2).Value = 96
$ws.Cells.Item(797, 3).Value = 43213
$ws.Cells.Item(797, 4).Value = 'April'
$ws.Cells.Item(797, 5).Value = 2018
$ws.Cells.Item(797, 6).Value = 'Monday'
$ws.Cells.Item(797, 7).Value = 'Leg Extension'
$ws.Cells.Item(797, 8).Value = 108
$ws.Cells.Item(797, 9).Value = 4
$ws.Cells.Item(797, 10).Value = 8
$ws.Cells.Item(797, 11).Value = 'Legs'

$ws.Cells.Item(798, 1).Value = 797
$ws.Cells.Item(798, 2).Value = 96
$ws.Cells.Item(798, 3).Value = 43213
$ws.Cells.Item(798, 4).Value = 'April'
$ws.Cells.Item(798, 5).Value = 2018
$ws.Cells.Item(798, 6).Value = 'Monday'
$ws.Cells.Item(798, 7).Value = 'Dumbell Rows'
$ws.Cells.Item(798, 8).Value = 30
$ws.Cells.Item(798, 9).Value = 4
$ws.Cells.Item(798, 10).Value = 8
$ws.Cells.Item(798, 11).Value = 'Back'

$ws.Cells.Item(799, 1).Value = 798
$ws.Cells.Item(799, 2).Value = 96
$ws.Cells.Item(799, 3).Value = 43213
$ws.Cells.Item(799, 4).Value = 'April'
$ws.Cells.Item(799, 5).Value = 2018
$ws.Cells.Item(799, 6).Value = 'Monday'
$ws.Cells.Item(799, 7).Value = 'Leg Raises'
$ws.Cells.Item(799, 8).Value = 0
$ws.Cells.Item(799, 9).Value = 4
$ws.Cells.Item(799, 10).Value = 10
$ws.Cells.Item(799, 11).Value = 'Core'

$ws.Cells.Item(800, 1).Value = 799
$ws.Cells.Item(800, 2).Value = 96
$ws.Cells.Item(800, 3).Value = 43213
$ws.Cells.Item(800, 4).Value = 'April'
$ws.Cells.Item(800, 5).Value = 2018
$ws.Cells.Item(800, 6).Value = 'Monday'
$ws.Cells.Item(800, 7).Value = 'Plank'
$ws.Cells.Item(800, 8).Value = 0
$ws.Cells.Item(800, 9).Value = 4
$ws.Cells.Item(800, 10).Value = 30
$ws.Cells.Item(800, 11).Value = 'Core'

$ws.Cells.Item(801, 1).Value = 800
$ws.Cells.Item(801, 2).Value = 97
$ws.Cells.Item(801, 3).Value = 43215
$ws.Cells.Item(801, 4).Value = 'April'
$ws.Cells.Item(801, 5).Value = 2018
$ws.Cells.Item(801, 6).Value = 'Wednesday'
$ws.Cells.Item(801, 7).Value = 'Pec Fly'
$ws.Cells.Item(801, 8).Value = 110
$ws.Cells.Item(801, 9).Value = 4
$ws.Cells.Item(801, 10).Value = 8
$ws.Cells.Item(801, 11).Value = 'Chest'

$ws.Cells.Item(802, 1).Value = 801
$ws.Cells.Item(802, 2).Value = 97
$ws.Cells.Item(802, 3).Value = 43215
$ws.Cells.Item(802, 4).Value = 'April'
$ws.Cells.Item(802, 5).Value = 2018
$ws.Cells.Item(802, 6).Value = 'Wednesday'
$ws.Cells.Item(802, 7).Value = 'Bicep Curl'
$ws.Cells.Item(802, 8).Value = 35
$ws.Cells.Item(802, 9).Value = 4
$ws.Cells.Item(802, 10).Value = 8
$ws.Cells.Item(802, 11).Value = 'Arms'

$ws.Cells.Item(803, 1).Value = 802
$ws.Cells.Item(803, 2).Value = 97
$ws.Cells.Item(803, 3).Value = 43215
$ws.Cells.Item(803, 4).Value = 'April'
$ws.Cells.Item(803, 5).Value = 2018
$ws.Cells.Item(803, 6).Value = 'Wednesday'
$ws.Cells.Item(803, 7).Value = 'Row (machine)'
$ws.Cells.Item(803, 8).Value = 60
$ws.Cells.Item(803, 9).Value = 4
$ws.Cells.Item(803, 10).Value = 8
$ws.Cells.Item(803, 11).Value = 'Back'

$ws.Cells.Item(804, 1).Value = 803
$ws.Cells.Item(804, 2).Value = 97
$ws.Cells.Item(804, 3).Value = 43215
$ws.Cells.Item(804, 4).Value = 'April'
$ws.Cells.Item(804, 5).Value = 2018
$ws.Cells.Item(804, 6).Value = 'Wednesday'
$ws.Cells.Item(804, 7).Value = 'Russian Twists'
$ws.Cells.Item(804, 8).Value = 15
$ws.Cells.Item(804, 9).Value = 4
$ws.Cells.Item(804, 10).Value = 12
$ws.Cells.Item(804, 11).Value = 'Core'

$ws.Cells.Item(805, 1).Value = 804
$ws.Cells.Item(805, 2).Value = 97
$ws.Cells.Item(805, 3).Value = 43215
$ws.Cells.Item(805, 4).Value = 'April'
$ws.Cells.Item(805, 5).Value = 2018
$ws.Cells.Item(805, 6).Value = 'Wednesday'
$ws.Cells.Item(805, 7).Value = 'Heel-taps'
$ws.Cells.Item(805, 8).Value = 0
$ws.Cells.Item(805, 9).Value = 3
$ws.Cells.Item(805, 10).Value = 12
$ws.Cells.Item(805, 11).Value = 'Core'

$ws.Cells.Item(806, 1).Value = 805
$ws.Cells.Item(806, 2).Value = 97
$ws.Cells.Item(806, 3).Value = 43215
$ws.Cells.Item(806, 4).Value = 'April'
$ws.Cells.Item(806, 5).Value = 2018
$ws.Cells.Item(806, 6).Value = 'Wednesday'
$ws.Cells.Item(806, 7).Value = 'Raised leg circles'
$ws.Cells.Item(806, 8).Value = 0
$ws.Cells.Item(806, 9).Value = 3
$ws.Cells.Item(806, 10).Value = 10
$ws.Cells.Item(806, 11).Value = 'Core'

$ws.Cells.Item(807, 1).Value = 806
$ws.Cells.Item(807, 2).Value = 97
$ws.Cells.Item(807, 3).Value = 43215
$ws.Cells.Item(807, 4).Value = 'April'
$ws.Cells.Item(807, 5).Value = 2018
$ws.Cells.Item(807, 6).Value = 'Wednesday'
$ws.Cells.Item(807, 7).Value = 'Scissors'
$ws.Cells.Item(807, 8).Value = 0
$ws.Cells.Item(807, 9).Value = 3
$ws.Cells.Item(807, 10).Value = 12
$ws.Cells.Item(807, 11).Value = 'Core'

$ws.Cells.Item(808, 1).Value = 807
$ws.Cells.Item(808, 2).Value = 97
$ws.Cells.Item(808, 3).Value = 43215
$ws.Cells.Item(808, 4).Value = 'April'
$ws.Cells.Item(808, 5).Value = 2018
$ws.Cells.Item(808, 6).Value = 'Wednesday'
$ws.Cells.Item(808, 7).Value = 'Knee-Pull ins'
$ws.Cells.Item(808, 8).Value = 0
$ws.Cells.Item(808, 9).Value = 3
$ws.Cells.Item(808, 10).Value = 10
$ws.Cells.Item(808, 11).Value = 'Core'

$ws.Cells.Item(809, 1).Value = 808
$ws.Cells.Item(809, 2).Value = 97
$ws.Cells.Item(809, 3).Value = 43215
$ws.Cells.Item(809, 4).Value = 'April'
$ws.Cells.Item(809, 5).Value = 2018
$ws.Cells.Item(809, 6).Value = 'Wednesday'
$ws.Cells.Item(809, 7).Value = 'Flitter Kicks'
$ws.Cells.Item(809, 8).Value = 0
$ws.Cells.Item(809, 9).Value = 3
$ws.Cells.Item(809, 10).Value = 20
$ws.Cells.Item(809, 11).Value = 'Core'

$ws.Cells.Item(810, 1).Value = 809
$ws.Cells.Item(810, 2).Value = 98
$ws.Cells.Item(810, 3).Value = 43216
$ws.Cells.Item(810, 4).Value = 'April'
$ws.Cells.Item(810, 5).Value = 2018
$ws.Cells.Item(810, 6).Value = 'Thursday'
$ws.Cells.Item(810, 7).Value = 'Bench Press'
$ws.Cells.Item(810, 8).Value = 85
$ws.Cells.Item(810, 9).Value = 5
$ws.Cells.Item(810, 10).Value = 5
$ws.Cells.Item(810, 11).Value = 'Chest'

$ws.Cells.Item(811, 1).Value = 810
$ws.Cells.Item(811, 2).Value = 98
$ws.Cells.Item(811, 3).Value = 43216
$ws.Cells.Item(811, 4).Value = 'April'
$ws.Cells.Item(811, 5).Value = 2018
$ws.Cells.Item(811, 6).Value = 'Thursday'
$ws.Cells.Item(811, 7).Value = 'Overhead Press'
$ws.Cells.Item(811, 8).Value = 50
$ws.Cells.Item(811, 9).Value = 5
$ws.Cells.Item(811, 10).Value = 5
$ws.Cells.Item(811, 11).Value = 'Shoulders'

$ws.Cells.Item(812, 1).Value = 811
$ws.Cells.Item(812, 2).Value = 98
$ws.Cells.Item(812, 3).Value = 43216
$ws.Cells.Item(812, 4).Value = 'April'
$ws.Cells.Item(812, 5).Value = 2018
$ws.Cells.Item(812, 6).Value = 'Thursday'
$ws.Cells.Item(812, 7).Value = 'Barbell Row'
$ws.Cells.Item(812, 8).Value = 100
$ws.Cells.Item(812, 9).Value = 5
$ws.Cells.Item(812, 10).Value = 5
$ws.Cells.Item(812, 11).Value = 'Back'

$ws.Cells.Item(813, 1).Value = 812
$ws.Cells.Item(813, 2).Value = 98
$ws.Cells.Item(813, 3).Value = 43216
$ws.Cells.Item(813, 4).Value = 'April'
$ws.Cells.Item(813, 5).Value = 2018
$ws.Cells.Item(813, 6).Value = 'Thursday'
$ws.Cells.Item(813, 7).Value = 'Upright Rows'
$ws.Cells.Item(813, 8).Value = 30
$ws.Cells.Item(813, 9).Value = 4
$ws.Cells.Item(813, 10).Value = 12
$ws.Cells.Item(813, 11).Value = 'Shoulders'

$ws.Cells.Item(814, 1).Value = 813
$ws.Cells.Item(814, 2).Value = 98
$ws.Cells.Item(814, 3).Value = 43216
$ws.Cells.Item(814, 4).Value = 'April'
$ws.Cells.Item(814, 5).Value = 2018
$ws.Cells.Item(814, 6).Value = 'Thursday'
$ws.Cells.Item(814, 7).Value = 'Press-up hold'
$ws.Cells.Item(814, 8).Value = 0
$ws.Cells.Item(814, 9).Value = 1
$ws.Cells.Item(814, 10).Value = 60
$ws.Cells.Item(814, 11).Value = 'Shoulders'

$ws.Cells.Item(815, 1).Value = 814
$ws.Cells.Item(815, 2).Value = 99
$ws.Cells.Item(815, 3).Value = 43219
$ws.Cells.Item(815, 4).Value = 'April'
$ws.Cells.Item(815, 5).Value = 2018
$ws.Cells.Item(815, 6).Value = 'Sunday'
$ws.Cells.Item(815, 7).Value = 'Deadlift'
$ws.Cells.Item(815, 8).Value = 130
$ws.Cells.Item(815, 9).Value = 5
$ws.Cells.Item(815, 10).Value = 3
$ws.Cells.Item(815, 11).Value = 'Legs'

$ws.Cells.Item(816, 1).Value = 815
$ws.Cells.Item(816, 2).Value = 99
$ws.Cells.Item(816, 3).Value = 43219
$ws.Cells.Item(816, 4).Value = 'April'
$ws.Cells.Item(816, 5).Value = 2018
$ws.Cells.Item(816, 6).Value = 'Sunday'
$ws.Cells.Item(816, 7).Value = 'Barbell Lunge'
$ws.Cells.Item(816, 8).Value = 40
$ws.Cells.Item(816, 9).Value = 3
$ws.Cells.Item(816, 10).Value = 8
$ws.Cells.Item(816, 11).Value = 'Legs'

$ws.Cells.Item(817, 1).Value = 816
$ws.Cells.Item(817, 2).Value = 99
$ws.Cells.Item(817, 3).Value = 43219
$ws.Cells.Item(817, 4).Value = 'April'
$ws.Cells.Item(817, 5).Value = 2018
$ws.Cells.Item(817, 6).Value = 'Sunday'
$ws.Cells.Item(817, 7).Value = 'Barbell Squat'
$ws.Cells.Item(817, 8).Value = 70
$ws.Cells.Item(817, 9).Value = 3
$ws.Cells.Item(817, 10).Value = 8
$ws.Cells.Item(817, 11).Value = 'Legs'

$ws.Cells.Item(818, 1).Value = 817
$ws.Cells.Item(818, 2).Value = 99
$ws.Cells.Item(818, 3).Value = 43219
$ws.Cells.Item(818, 4).Value = 'April'
$ws.Cells.Item(818, 5).Value = 2018
$ws.Cells.Item(818, 6).Value = 'Sunday'
$ws.Cells.Item(818, 7).Value = 'Leg Press'
$ws.Cells.Item(818, 8).Value = 120
$ws.Cells.Item(818, 9).Value = 3
$ws.Cells.Item(818, 10).Value = 6
$ws.Cells.Item(818, 11).Value = 'Legs'

$ws.Range("A819").Select() | Out-Null
